$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 47
$ws.Range("H47").Value = 8000
$ws.Range("I47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("M47").ClearContents()
# Row 113
$ws.Range("J113").Value = 5138.1
$ws.Range("L113").Value = 5138.1
$ws.Range("N113").Value = -11646.1
# Row 132
$ws.Range("H132").Value = 29414616
$ws.Range("I132").Value = 31252952
$ws.Range("J132").Value = 1250
$ws.Range("K132").Value = 93758856
$ws.Range("L132").Value = 3750
$ws.Range("M132").Value = -93756326
$ws.Range("N132").Value = -8810
# Row 137
$ws.Range("H137").Value = 92545.42999999999
$ws.Range("I137").Value = 134704.67
$ws.Range("J137").Value = 2204.2144
$ws.Range("K137").Value = 404114.01
$ws.Range("L137").Value = 6612.6432
$ws.Range("M137").Value = -401564.01
$ws.Range("N137").Value = -11712.6432

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 5097.22
$ws.Range("I32").Value = 4811.567
$ws.Range("J32").Value = 14333.333
$ws.Range("K32").Value = 4811.567
$ws.Range("L32").Value = 14333.333
$ws.Range("M32").Value = -4524.567
$ws.Range("N32").Value = -14907.333
# Row 61
$ws.Range("H61").Value = 9262441
$ws.Range("I61").Value = 11114026
$ws.Range("J61").Value = 4515.6665
$ws.Range("K61").Value = 11114026
$ws.Range("L61").Value = 4515.6665
$ws.Range("M61").Value = -11113814
$ws.Range("N61").Value = -4939.6665
# Row 74
$ws.Range("H74").Value = 22728562
$ws.Range("I74").Value = 32258668
$ws.Range("K74").Value = 32258668
$ws.Range("M74").Value = -32257794
# Row 77
$ws.Range("H77").Value = 22728562
$ws.Range("I77").Value = 32258668
$ws.Range("K77").Value = 161293340
$ws.Range("M77").Value = -161288972
# Row 132
$ws.Range("H132").Value = 13905410
$ws.Range("I132").Value = 17859530
$ws.Range("J132").Value = 65994.25
$ws.Range("K132").Value = 53578590
$ws.Range("L132").Value = 197982.75
$ws.Range("M132").Value = -53576060
$ws.Range("N132").Value = -203042.75
# Row 136
$ws.Range("H136").Value = 9262441
$ws.Range("I136").Value = 11114026
$ws.Range("J136").Value = 4515.6665
$ws.Range("K136").Value = 33342078
$ws.Range("L136").Value = 13546.9995
$ws.Range("M136").Value = -33339528
$ws.Range("N136").Value = -18646.9995

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 80
$ws.Range("H80").Value = 1001.2143
$ws.Range("J80").Value = 1187.6666
$ws.Range("L80").Value = 1187.6666
$ws.Range("N80").Value = -3183.6666
# Row 83
$ws.Range("H83").Value = 1001.2143
$ws.Range("J83").Value = 1187.6666
$ws.Range("L83").Value = 5938.333000000001
$ws.Range("N83").Value = -15922.333

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 6269.2666
$ws.Range("I31").Value = 4384.4443
$ws.Range("J31").Value = 7077.048
$ws.Range("K31").Value = 4384.4443
$ws.Range("L31").Value = 7077.048
$ws.Range("M31").Value = -4089.4443
$ws.Range("N31").Value = -7667.048
# Row 34
$ws.Range("H34").Value = 6269.2666
$ws.Range("I34").Value = 4384.4443
$ws.Range("J34").Value = 7077.048
$ws.Range("K34").Value = 4384.4443
$ws.Range("L34").Value = 7077.048
$ws.Range("M34").Value = -4182.4443
$ws.Range("N34").Value = -7481.048
# Row 52
$ws.Range("H52").Value = 25157.777
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 25157.777
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 25157.777
$ws.Range("M52").ClearContents()
$ws.Range("N52").Value = -25745.777
# Row 58
$ws.Range("H58").Value = 14155.462
$ws.Range("I58").Value = 1315.7727
$ws.Range("J58").Value = 30771.53
$ws.Range("K58").Value = 1315.7727
$ws.Range("L58").Value = 30771.53
$ws.Range("M58").Value = -1112.7727
$ws.Range("N58").Value = -31177.53
# Row 86
$ws.Range("H86").Value = 8171.8887
$ws.Range("I86").Value = 1452
$ws.Range("J86").Value = 16571.75
$ws.Range("K86").Value = 1452
$ws.Range("L86").Value = 16571.75
$ws.Range("M86").Value = -329
$ws.Range("N86").Value = -18817.75
# Row 89
$ws.Range("H89").Value = 8171.8887
$ws.Range("I89").Value = 1452
$ws.Range("J89").Value = 16571.75
$ws.Range("K89").Value = 7260
$ws.Range("L89").Value = 82858.75
$ws.Range("M89").Value = -1644
$ws.Range("N89").Value = -94090.75
# Row 132
$ws.Range("H132").Value = 43481250
$ws.Range("I132").Value = 55557520
$ws.Range("J132").Value = 6662.6
$ws.Range("K132").Value = 166672560
$ws.Range("L132").Value = 19987.8
$ws.Range("M132").Value = -166670030
$ws.Range("N132").Value = -25047.8
# Row 134
$ws.Range("H134").Value = 58824350
$ws.Range("I134").Value = 58824350
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 176473050
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -176470515
$ws.Range("N134").ClearContents()
# Row 136
$ws.Range("H136").Value = 14155.462
$ws.Range("I136").Value = 1315.7727
$ws.Range("J136").Value = 30771.53
$ws.Range("K136").Value = 3947.3181
$ws.Range("L136").Value = 92314.59
$ws.Range("M136").Value = -1397.3181
$ws.Range("N136").Value = -97414.59

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1294.119
$ws.Range("I5").Value = 870.5925999999999
$ws.Range("K5").Value = 2611.7778
$ws.Range("M5").Value = -2499.7778
# Row 131
$ws.Range("H131").Value = 695.71
$ws.Range("J131").Value = 737.75
$ws.Range("L131").Value = 2213.25
$ws.Range("N131").Value = -12293.25
# Row 135
$ws.Range("H135").Value = 1294.119
$ws.Range("I135").Value = 870.5925999999999
$ws.Range("K135").Value = 7835.3334
$ws.Range("M135").Value = -5300.3334
# Row 141
$ws.Range("H141").Value = 4805.8
$ws.Range("I141").Value = 4805.8
$ws.Range("K141").Value = 14417.4
$ws.Range("M141").Value = -9237.400000000001

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 126
$ws.Range("H126").Value = 3879.5789
$ws.Range("I126").Value = 2742.6667
$ws.Range("J126").Value = 5828.5713
$ws.Range("K126").Value = 8228.000100000001
$ws.Range("L126").Value = 17485.7139
$ws.Range("M126").Value = -5758.000100000001
$ws.Range("N126").Value = -22425.7139
# Row 132
$ws.Range("H132").Value = 5104526.5
$ws.Range("I132").Value = 6689008.5
$ws.Range("J132").Value = 86999.5
$ws.Range("K132").Value = 20067025.5
$ws.Range("L132").Value = 260998.5
$ws.Range("M132").Value = -20064495.5
$ws.Range("N132").Value = -266058.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 69
$ws.Range("H69").Value = 30000
$ws.Range("J69").Value = 30000
$ws.Range("L69").Value = 30000
$ws.Range("N69").Value = -31622
# Row 72
$ws.Range("H72").Value = 30000
$ws.Range("J72").Value = 30000
$ws.Range("L72").Value = 90000
$ws.Range("N72").Value = -98112
# Row 136
$ws.Range("H136").Value = 2451.6897
$ws.Range("I136").Value = 2451.6897
$ws.Range("K136").Value = 7355.0691
$ws.Range("M136").Value = -4805.0691

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 25001346
$ws.Range("I132").Value = 27778718
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 83336154
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -83333624
$ws.Range("N132").Value = -20057
# Row 136
$ws.Range("H136").Value = 27780868
$ws.Range("I136").Value = 35715856
$ws.Range("J136").Value = 8413.125
$ws.Range("K136").Value = 107147568
$ws.Range("L136").Value = 25239.375
$ws.Range("M136").Value = -107145018
$ws.Range("N136").Value = -30339.375
